# Apply the "pushing new eda2 & processing2files" edit:
#  - Recode the "Hair Color" values on the Data sheet from full color
#    names to short codes (Blonde->BL, Red->R, Black->BL, Brown->BR,
#    White->W, Purple->P).
#  - Add two new rows to the Codebook sheet describing "Strands of Hair"
#    and "Hair Color" (with the new code legend).
#  - Update the active sheet / selection so the Codebook sheet ends up
#    selected (matches the saved view state in the target workbook).

$wb = $excel.ActiveWorkbook

$dataWs = $wb.Worksheets.Item("Data")
$codebookWs = $wb.Worksheets.Item("Codebook")

# --- Data sheet: recode Hair Color column (E) -----------------------------
# Full color names are replaced by short codes everywhere:
#   Blonde -> BL, Red -> R, Black -> BL, Brown -> BR, White -> W, Purple -> P
$dataWs.Range("E2").Value = "BL"   # was Blonde
$dataWs.Range("E3").Value = "R"    # was Red
$dataWs.Range("E4").Value = "BL"   # was Black
$dataWs.Range("E5").Value = "BR"   # was Brown
$dataWs.Range("E6").Value = "W"    # was White
$dataWs.Range("E7").Value = "P"    # was Purple
$dataWs.Range("E8").Value = "R"    # was Red
$dataWs.Range("E9").Value = "BR"   # was Brown
$dataWs.Range("E10").Value = "BR"  # was Brown
$dataWs.Range("E11").Value = "BR"  # was Brown
$dataWs.Range("E12").Value = "BL"  # was Blonde
$dataWs.Range("E13").Value = "W"   # was White
$dataWs.Range("E14").Value = "BL"  # was Black
$dataWs.Range("E15").Value = "BR"  # was Brown

# --- Codebook sheet: document the two hair-related variables --------------
$codebookWs.Range("A5").Value = "Strands of Hair"
$codebookWs.Range("B5").Value = "How many strands of hair the individual has"
$codebookWs.Range("C5").Value = "numeric value >0 or NA"

$codebookWs.Range("A6").Value = "Hair Color"
$codebookWs.Range("B6").Value = "What color the invididual's hair is"
$codebookWs.Range("C6").Value = "BR/BL/R/W/P; BR=brown, BL= black, R= red, W= white, P= purple"

# --- Column widths on Codebook (minor width nudge from the source file) ---
$codebookWs.Columns.Item(1).ColumnWidth = 13.498697916666666
$codebookWs.Columns.Item(2).ColumnWidth = 29.498697916666668
$codebookWs.Columns.Item(3).ColumnWidth = 20.998697916666668

# --- View state: selections + which sheet/cell is active ------------------
$dataWs.Range("E5").Select() | Out-Null
$codebookWs.Range("C9").Select() | Out-Null
$codebookWs.Activate() | Out-Null
